$d = $word.ActiveDocument

# Insert two new paragraphs after the existing content:
#  1) an empty paragraph (just a paragraph mark carrying the en-US lang rPr)
#  2) a paragraph containing the new sentence
# Both are appended in a single InsertXML call at the very end of the
# document's content so the engine inserts them as genuine new paragraphs
# rather than merging into / replacing the trailing paragraph mark.

$w_ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$emptyParaXml = '<w:p ' + $w_ns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$textParaXml  = '<w:p ' + $w_ns + '><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Saya adalah seorang pelajar di sebuah universitas</w:t></w:r></w:p>'

$insertRange = $d.Range($d.Content.End, $d.Content.End)
$insertRange.InsertXML($emptyParaXml + $textParaXml)
